$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.450441837310791
$ws.Range("B1").Value = 1.934596061706543
$ws.Range("C1").Value = 3.098910570144653
$ws.Range("D1").Value = 4.300515174865723
$ws.Range("E1").Value = 1.015467524528503
